# Add 7-31-19 Results for Subchallenge 2
# Populates column D ("7/31 Submission") on the SubCh2 worksheet with the
# new round of vote results, which also updates the dependent SUM formulas
# in column E and the correlation figures in I5:J5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SubCh2")

# New "7/31 Submission" (column D) values for each isolate row (2-33),
# in row order.
$dValues = @(0, 1, 1, 1, 1, 0, 1, 0, 1, 0, 0, 1, 1, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 1, 0, 1, 1, 0, 0, 1, 0, 1)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Update the pairwise correlation/vote-agreement figures now that the
# 7/31 Submission column has real data (previously #DIV/0! with no data).
$ws.Range("I5").Value = 0.46666666666666667
$ws.Range("J5").Value = 0.57594696646956689

# Highlight the Vote column with a 3-colour scale (green low -> red high).
$voteRange = $ws.Range("E2:E33")
$colorScale = $voteRange.FormatConditions.AddColorScale(3)
$colorScale.ColorScaleCriteria.Item(1).FormatColor.Color = 8109667   # FF63BE7B green
$colorScale.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167   # FFFFEB84 yellow
$colorScale.ColorScaleCriteria.Item(3).FormatColor.Color = 7039480   # FFF8696B red

# Reflect where the author was last working on this sheet.
$ws.Activate()
$ws.Range("F7").Select()
